$d = $word.ActiveDocument

$nbsp = [char]0x00A0
$apos = [char]0x2019
$pkgOpen  = "<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'><pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'><pkg:xmlData><w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:body>"
$pkgClose = "</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"

# ------------------------------------------------------------------
# 1) Insert three new blank paragraphs right before the paragraph
#    that currently holds only the _GoBack bookmark (paragraph 6).
# ------------------------------------------------------------------
$bookmarkPara = $d.Paragraphs.Item(6)
$insertPoint = $bookmarkPara.Range.Duplicate
$insertPoint.Collapse(1)
$insertPoint.InsertParagraphBefore()
$insertPoint.InsertParagraphBefore()
$insertPoint.InsertParagraphBefore()

# ------------------------------------------------------------------
# 2) The bookmark paragraph is now paragraph 9. Drop the _GoBack
#    bookmark from it and give it the "Q2 : " text instead.
# ------------------------------------------------------------------
$q2Para = $d.Paragraphs.Item(9)
$d.Bookmarks.Item("_GoBack").Delete()

$q2Range = $d.Range($q2Para.Range.Start, $q2Para.Range.End)
$q2Xml = $pkgOpen + "<w:p><w:r><w:t xml:space='preserve'>Q2" + $nbsp + ": </w:t></w:r></w:p>" + $pkgClose
[void]$q2Range.InsertXML($q2Xml)

# ------------------------------------------------------------------
# 3) The paragraph that used to hold "Q2 : " is now paragraph 12.
#    Replace its text with the new annotation sentence (with the
#    spell-check marks around "l'arduino", exactly as in the source).
#    A trailing placeholder character is appended so that the
#    _GoBack bookmark can be anchored immediately after the real
#    text (collapsed ranges sitting exactly on a paragraph-end
#    boundary are not placed reliably), then it is deleted again.
# ------------------------------------------------------------------
$annotationPara = $d.Paragraphs.Item(12)
$annotationRange = $d.Range($annotationPara.Range.Start, $annotationPara.Range.End)
$annotationXml = $pkgOpen + "<w:p>" +
    "<w:r><w:t xml:space='preserve'>PASSER LE AREF ET LE VCC </w:t></w:r>" +
    "<w:r><w:t xml:space='preserve'>de la </w:t></w:r>" +
    "<w:proofErr w:type='spellStart'/>" +
    "<w:r><w:t>l" + $apos + "arduino</w:t></w:r>" +
    "<w:proofErr w:type='spellEnd'/>" +
    "<w:r><w:t xml:space='preserve'> à 5V#</w:t></w:r>" +
    "</w:p>" + $pkgClose
$annotationRange.InsertXML($annotationXml)

$docEnd = $d.Content.End
$placeholder = $d.Range($docEnd - 2, $docEnd - 1)
$d.Bookmarks.Add("_GoBack", $placeholder)
$placeholder = $d.Range($docEnd - 2, $docEnd - 1)
$placeholder.Text = ""
